$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values remain stored as text, matching the
# original inlineStr cell type (values like "1.00" or "0.0000227" must
# not be auto-converted to numbers which would drop formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.797.86'
$ws.Range("E2").Value = '  +4.63%  '

$ws.Range("D3").Value = '3.619.12'
$ws.Range("E3").Value = '  +4.38%  '

$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '632.16'
$ws.Range("E5").Value = '  +4.76%  '

$ws.Range("D6").Value = '159.50'
$ws.Range("E6").Value = '  +8.06%  '

$ws.Range("D7").Value = '3.619.17'
$ws.Range("E7").Value = '  +4.27%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = '0.496'
$ws.Range("E9").Value = '  +3.74%  '

$ws.Range("D10").Value = '0.150'
$ws.Range("E10").Value = '  +9.87%  '

$ws.Range("D11").Value = '7.48'
$ws.Range("E11").Value = '  +8.56%  '

$ws.Range("D12").Value = '0.444'
$ws.Range("E12").Value = '  +5.87%  '

$ws.Range("D13").Value = '0.0000227'
$ws.Range("E13").Value = '  +5.75%  '

$ws.Range("D14").Value = '33.95'
$ws.Range("E14").Value = '  +9.26%  '

$ws.Range("D15").Value = '4.220.02'
$ws.Range("E15").Value = '  +4.20%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '70.006.69'
$ws.Range("E16").Value = '  +5.07%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.613.26'
$ws.Range("E17").Value = '  +4.94%  '

$ws.Range("E18").Value = '  +1.18%  '

$ws.Range("D19").Value = '6.79'
$ws.Range("E19").Value = '  +5.96%  '

$ws.Range("D20").Value = '16.21'
$ws.Range("E20").Value = '  +8.77%  '

$ws.Range("D21").Value = '10.25'
$ws.Range("E21").Value = '  +14.22%  '

$ws.Range("D22").Value = '465.91'
$ws.Range("E22").Value = '  +5.74%  '

$ws.Range("D23").Value = '0.649'
$ws.Range("E23").Value = '  +4.51%  '

$ws.Range("D24").Value = '79.08'
$ws.Range("E24").Value = '  +3.01%  '

$ws.Range("D25").Value = '0.0000137'
$ws.Range("E25").Value = '  +10.13%  '

$ws.Range("D26").Value = '10.82'
$ws.Range("E26").Value = '  +7.79%  '

$ws.Range("D27").Value = '3.759.02'
$ws.Range("E27").Value = '  +4.37%  '

$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("D29").Value = '9.42'
$ws.Range("E29").Value = '  +13.97%  '

$ws.Range("D30").Value = '2.67'
$ws.Range("E30").Value = '  +6.27%  '

$ws.Range("D31").Value = '1.75'
$ws.Range("E31").Value = '  +13.77%  '

$ws.Range("D32").Value = '0.175'
$ws.Range("E32").Value = '  +9.70%  '

$ws.Range("D33").Value = '6.62'
$ws.Range("E33").Value = '  +8.37%  '

$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("D35").Value = '1.98'
$ws.Range("E35").Value = '  +7.01%  '

$ws.Range("D36").Value = '26.68'
$ws.Range("E36").Value = '  +4.65%  '

$ws.Range("D37").Value = '3.611.97'
$ws.Range("E37").Value = '  +4.74%  '

$ws.Range("D38").Value = '8.52'
$ws.Range("E38").Value = '  +7.76%  '

$ws.Range("D39").Value = '2.46'
$ws.Range("E39").Value = '  +15.48%  '

$ws.Range("E40").Value = '  +0.04%  '

$ws.Range("D41").Value = '0.0930'
$ws.Range("E41").Value = '  +7.62%  '

$ws.Range("D42").Value = '179.05'
$ws.Range("E42").Value = '  +3.71%  '

$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("D44").Value = '5.73'
$ws.Range("E44").Value = '  +5.90%  '

$ws.Range("D45").Value = '31.98'
$ws.Range("E45").Value = '  +23.89%  '

$ws.Range("D46").Value = '0.916'
$ws.Range("E46").Value = '  +4.44%  '

$ws.Range("D47").Value = '1.40'

$ws.Range("D48").Value = '2.80'
$ws.Range("E48").Value = '  +12.93%  '

$ws.Range("D49").Value = '46.11'
$ws.Range("E49").Value = '  +1.39%  '

$ws.Range("D50").Value = '7.87'
$ws.Range("E50").Value = '  +4.32%  '

$ws.Range("D51").Value = '0.271'
$ws.Range("E51").Value = '  +11.18%  '
